$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5971311412731097
$ws.Range("C2").Value = 0.1924353492518591
$ws.Range("D2").Value = 0.0526970135717022
$ws.Range("E2").Value = 0.1283788784727697
$ws.Range("F2").Value = 1.056645474702407
$ws.Range("K2").Value = 0.297797895826875
$ws.Range("L2").Value = 0.1908954316181024
$ws.Range("M2").Value = 0.1672341601680216
$ws.Range("N2").Value = 2.058920891053175
$ws.Range("O2").Value = 3.80778335846091

$ws.Range("B3").Value = 0.5617636639495345
$ws.Range("C3").Value = 0.1918725758464177
$ws.Range("D3").Value = 0.05072020776029262
$ws.Range("E3").Value = 0.1288173709435192
$ws.Range("F3").Value = 1.056773709657563
$ws.Range("K3").Value = 0.2653995575084167
$ws.Range("L3").Value = 0.1883031288397063
$ws.Range("M3").Value = 0.1606412002114723
$ws.Range("N3").Value = 2.078227358966034
$ws.Range("O3").Value = 3.821891200383845

$ws.Range("B4").Value = 0.5402721185319592
$ws.Range("C4").Value = 0.1915311518132796
$ws.Range("D4").Value = 0.04949122407847995
$ws.Range("E4").Value = 0.1291335681084753
$ws.Range("F4").Value = 1.057390381255011
$ws.Range("K4").Value = 0.2455342934939893
$ws.Range("L4").Value = 0.1868017948968728
$ws.Range("M4").Value = 0.156667128646923
$ws.Range("N4").Value = 2.090690007394027
$ws.Range("O4").Value = 3.83233153525569

$ws.Range("B5").Value = 0.5315711090181878
$ws.Range("C5").Value = 0.1913930744892447
$ws.Range("D5").Value = 0.04898659570397257
$ws.Range("E5").Value = 0.1292742507469953
$ws.Range("F5").Value = 1.057777067769408
$ws.Range("K5").Value = 0.2374463671372808
$ws.Range("L5").Value = 0.1862127748914162
$ws.Range("M5").Value = 0.155066396244095
$ws.Range("N5").Value = 2.09592165782551
$ws.Range("O5").Value = 3.837033474484059

$ws.Range("B6").Value = 0.530129771082926
$ws.Range("C6").Value = 0.1913702111918987
$ws.Range("D6").Value = 0.04890257313588364
$ws.Range("E6").Value = 0.1292983260959915
$ws.Range("F6").Value = 1.057849457388762
$ws.Range("K6").Value = 0.2361038284505668
$ws.Range("L6").Value = 0.1861163467261875
$ws.Range("M6").Value = 0.1548017308078187
$ws.Range("N6").Value = 2.096799612436003
$ws.Range("O6").Value = 3.837841259993681

$ws.Range("B7").Value = 0.5401545422300273
$ws.Range("C7").Value = 0.191529285351745
$ws.Range("D7").Value = 0.04948443387566215
$ws.Range("E7").Value = 0.1291354174775332
$ws.Range("F7").Value = 1.057395047905985
$ws.Range("K7").Value = 0.2454251866438568
$ws.Range("L7").Value = 0.1867937588201869
$ws.Range("M7").Value = 0.1566454646000999
$ws.Range("N7").Value = 2.090759943645171
$ws.Range("O7").Value = 3.832393135265278

$ws.Range("B8").Value = 0.5848902359443571
$ws.Range("C8").Value = 0.1922404610958424
$ws.Range("D8").Value = 0.05201858002074289
$ws.Range("E8").Value = 0.1285203363893306
$ws.Range("F8").Value = 1.056578110902407
$ws.Range("K8").Value = 0.2866214993042036
$ws.Range("L8").Value = 0.1899828919965003
$ws.Range("M8").Value = 0.1649456105346871
$ws.Range("N8").Value = 2.065451509691917
$ws.Range("O8").Value = 3.812278885906579

$ws.Range("B9").Value = 0.6743749815216802
$ws.Range("C9").Value = 0.193666997782401
$ws.Range("D9").Value = 0.05686666414143104
$ws.Range("E9").Value = 0.1276859123110565
$ws.Range("F9").Value = 1.059240366617011
$ws.Range("K9").Value = 0.3676102577205995
$ws.Range("L9").Value = 0.1969512646504583
$ws.Range("M9").Value = 0.1818053654811038
$ws.Range("N9").Value = 2.020646882275122
$ws.Range("O9").Value = 3.786931994741565

$ws.Range("B10").Value = 0.7411701001496738
$ws.Range("C10").Value = 0.1947336084718572
$ws.Range("D10").Value = 0.06035401061168244
$ws.Range("E10").Value = 0.1272984062128035
$ws.Range("F10").Value = 1.063792412862043
$ws.Range("K10").Value = 0.4272222574185207
$ws.Range("L10").Value = 0.202504045389361
$ws.Range("M10").Value = 0.1945435030926106
$ws.Range("N10").Value = 1.990666286348555
$ws.Range("O10").Value = 3.776892552694818

$ws.Range("B11").Value = 0.7717809948944137
$ws.Range("C11").Value = 0.1952226733131539
$ws.Range("D11").Value = 0.06192420767465023
$ws.Range("E11").Value = 0.1271708716254984
$ws.Range("F11").Value = 1.066426347129152
$ws.Range("K11").Value = 0.454362425134633
$ws.Range("L11").Value = 0.2051237213648136
$ws.Range("M11").Value = 0.2004138683314736
$ws.Range("N11").Value = 1.977664267174015
$ws.Range("O11").Value = 3.774187001787482

$ws.Range("B12").Value = 0.783404461173177
$ws.Range("C12").Value = 0.1954084061647663
$ws.Range("D12").Value = 0.06251645292115882
$ws.Range("E12").Value = 0.1271295672909858
$ws.Range("F12").Value = 1.067504647034937
$ws.Range("K12").Value = 0.464642555866476
$ws.Range("L12").Value = 0.2061291401727061
$ws.Range("M12").Value = 0.2026476106202111
$ws.Range("N12").Value = 1.972832212803326
$ws.Range("O12").Value = 3.773429925600311

$ws.Range("B13").Value = 0.7808997362418211
$ws.Range("C13").Value = 0.1953683818153351
$ws.Range("D13").Value = 0.06238900738151898
$ws.Range("E13").Value = 0.1271381523288042
$ws.Range("F13").Value = 1.067268820010895
$ws.Range("K13").Value = 0.4624284316229819
$ws.Range("L13").Value = 0.2059120101633169
$ws.Range("M13").Value = 0.2021660576078688
$ws.Range("N13").Value = 1.973868810952698
$ws.Range("O13").Value = 3.773581083710127

$ws.Range("B14").Value = 0.7727366310934656
$ws.Range("C14").Value = 0.1952379430679514
$ws.Range("D14").Value = 0.06197297936176227
$ws.Range("E14").Value = 0.1271673335071526
$ws.Range("F14").Value = 1.066513439003174
$ws.Range("K14").Value = 0.4552081255426117
$ws.Range("L14").Value = 0.2052061695213752
$ws.Range("M14").Value = 0.2005974246849505
$ws.Range("N14").Value = 1.977264896621692
$ws.Range("O14").Value = 3.774119357669946

$ws.Range("B15").Value = 0.7677406139069376
$ws.Range("C15").Value = 0.1951581145924095
$ws.Range("D15").Value = 0.06171784291297655
$ws.Range("E15").Value = 0.1271861175924514
$ws.Range("F15").Value = 1.066061277133457
$ws.Range("K15").Value = 0.4507858219684522
$ws.Range("L15").Value = 0.2047755653740779
$ws.Range("M15").Value = 0.1996379895488047
$ws.Range("N15").Value = 1.979357017700412
$ws.Range("O15").Value = 3.774483890428712

$ws.Range("B16").Value = 0.7391740886671698
$ws.Range("C16").Value = 0.1947017228442718
$ws.Range("D16").Value = 0.06025106632855426
$ws.Range("E16").Value = 0.1273077201193136
$ws.Range("F16").Value = 1.063631603329782
$ws.Range("K16").Value = 0.425448992883986
$ws.Range("L16").Value = 0.2023347228509493
$ws.Range("M16").Value = 0.1941613738473436
$ws.Range("N16").Value = 1.991528805310841
$ws.Range("O16").Value = 3.777106805329964

$ws.Range("B17").Value = 0.7217067441532663
$ws.Range("C17").Value = 0.1944227150797886
$ws.Range("D17").Value = 0.05934707739869083
$ws.Range("E17").Value = 0.1273947902348631
$ws.Range("F17").Value = 1.062285252098377
$ws.Range("K17").Value = 0.4099110581730088
$ws.Range("L17").Value = 0.2008612925735207
$ws.Range("M17").Value = 0.1908209519803279
$ws.Range("N17").Value = 1.999158794897243
$ws.Range("O17").Value = 3.779192485032894

$ws.Range("B18").Value = 0.7116812490530151
$ws.Range("C18").Value = 0.1942626015570355
$ws.Range("D18").Value = 0.05882560206748622
$ws.Range("E18").Value = 0.1274494598057263
$ws.Range("F18").Value = 1.061563882147468
$ws.Range("K18").Value = 0.4009761833239054
$ws.Range("L18").Value = 0.2000226368923563
$ws.Range("M18").Value = 0.1889067646117439
$ws.Range("N18").Value = 2.003607253207344
$ws.Range("O18").Value = 3.780567333431634

$ws.Range("B19").Value = 0.7082904593763999
$ws.Range("C19").Value = 0.1942084530649737
$ws.Range("D19").Value = 0.05864877828365422
$ws.Range("E19").Value = 0.1274687587540608
$ws.Range("F19").Value = 1.061328747182934
$ws.Range("K19").Value = 0.3979513693470267
$ws.Range("L19").Value = 0.1997401998166453
$ws.Range("M19").Value = 0.1882598829399882
$ws.Range("N19").Value = 2.005123710121826
$ws.Range("O19").Value = 3.781062934361728

$ws.Range("B20").Value = 0.7235639771267302
$ws.Range("C20").Value = 0.1944523783689505
$ws.Range("D20").Value = 0.0594434665392285
$ws.Range("E20").Value = 0.1273850466516659
$ws.Range("F20").Value = 1.062423087278844
$ws.Range("K20").Value = 0.4115648804274485
$ws.Range("L20").Value = 0.2010172291431758
$ws.Range("M20").Value = 0.1911758080665678
$ws.Range("N20").Value = 1.998340370879911
$ws.Range("O20").Value = 3.778952327753245

$ws.Range("B21").Value = 0.7751334753958758
$ws.Range("C21").Value = 0.1952762417412117
$ws.Range("D21").Value = 0.06209524095272201
$ws.Range("E21").Value = 0.127158572726227
$ws.Range("F21").Value = 1.066733118282983
$ws.Range("K21").Value = 0.4573288339641124
$ws.Range("L21").Value = 0.2054131287787868
$ws.Range("M21").Value = 0.2010578791196664
$ws.Range("N21").Value = 1.976264899653275
$ws.Range("O21").Value = 3.773953996568792

$ws.Range("B22").Value = 0.8090220938199764
$ws.Range("C22").Value = 0.1958177908865224
$ws.Range("D22").Value = 0.06381459554161495
$ws.Range("E22").Value = 0.1270512945969884
$ws.Range("F22").Value = 1.070021392605696
$ws.Range("K22").Value = 0.4872539609882551
$ws.Range("L22").Value = 0.208364202605793
$ws.Range("M22").Value = 0.2075790563325484
$ws.Range("N22").Value = 1.9623708139013
$ws.Range("O22").Value = 3.772246143260304

$ws.Range("B23").Value = 0.7909183860227529
$ws.Range("C23").Value = 0.1955284782941149
$ws.Range("D23").Value = 0.06289820775997157
$ws.Range("E23").Value = 0.1271048299292303
$ws.Range("F23").Value = 1.068223276733022
$ws.Range("K23").Value = 0.4712810758239812
$ws.Range("L23").Value = 0.206782035348084
$ws.Range("M23").Value = 0.2040928886248281
$ws.Range("N23").Value = 1.969737522654237
$ws.Range("O23").Value = 3.773015093655914

$ws.Range("B24").Value = 0.7227242699903798
$ws.Range("C24").Value = 0.1944389666863842
$ws.Range("D24").Value = 0.05939989448924621
$ws.Range("E24").Value = 0.1273894373591045
$ws.Range("F24").Value = 1.062360607937521
$ws.Range("K24").Value = 0.4108171932975608
$ws.Range("L24").Value = 0.2009467039288779
$ws.Range("M24").Value = 0.1910153581038827
$ws.Range("N24").Value = 1.998710187638434
$ws.Range("O24").Value = 3.779060355349259

$ws.Range("B25").Value = 0.649980945936079
$ws.Range("C25").Value = 0.1932777604034968
$ws.Range("D25").Value = 0.055568178970411
$ws.Range("E25").Value = 0.1278719599608173
$ws.Range("F25").Value = 1.058063999417158
$ws.Range("K25").Value = 0.3456803888757349
$ws.Range("L25").Value = 0.1949898972847279
$ws.Range("M25").Value = 0.1818053654811038
$ws.Range("N25").Value = 2.03225159848947
$ws.Range("O25").Value = 3.792280953399711

